$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "median" (C), "LB" (D), and "UB" (E) columns with new equation output.
# Each entry: row, median, LB (optional, only set when changed), UB

$ws.Range("C2").Value = 6
$ws.Range("E2").Value = 71

$ws.Range("C3").Value = 12
$ws.Range("E3").Value = 173

$ws.Range("C4").Value = 28
$ws.Range("E4").Value = 321

$ws.Range("C5").Value = 58
$ws.Range("E5").Value = 649.1

$ws.Range("C6").Value = 106.5
$ws.Range("E6").Value = 1412.7

$ws.Range("C7").Value = 230.5
$ws.Range("E7").Value = 2838.2

$ws.Range("C8").Value = 480
$ws.Range("E8").Value = 5635.3

$ws.Range("C9").Value = 875
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 12087.9

$ws.Range("C10").Value = 1857.5
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 23299.1

$ws.Range("C11").Value = 3263.5
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 50033.9

$ws.Range("C12").Value = 6520
$ws.Range("D12").Value = 7
$ws.Range("E12").Value = 92566.8

$ws.Range("C13").Value = 11888
$ws.Range("D13").Value = 23.9
$ws.Range("E13").Value = 165236.2

$ws.Range("C14").Value = 21216
$ws.Range("D14").Value = 19.9
$ws.Range("E14").Value = 298179.3

$ws.Range("C15").Value = 29629
$ws.Range("D15").Value = 83
$ws.Range("E15").Value = 359811.5

$ws.Range("C16").Value = 31906.5
$ws.Range("D16").Value = 39
$ws.Range("E16").Value = 379376.6

$ws.Range("C17").Value = 28060.5
$ws.Range("D17").Value = 52
$ws.Range("E17").Value = 315414.5

$ws.Range("C18").Value = 10134
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 198451.3

$ws.Range("C19").Value = 0
$ws.Range("E19").Value = 120622.9

$ws.Range("C20").Value = 0
$ws.Range("E20").Value = 47518

$ws.Range("C21").Value = 0
$ws.Range("E21").Value = 7756.5

$ws.Range("C22").Value = 0
$ws.Range("E22").Value = 0

$ws.Range("C23").Value = 0
$ws.Range("E23").Value = 0

$ws.Range("C24").Value = 0
$ws.Range("E24").Value = 0

$ws.Range("C25").Value = 0
$ws.Range("E25").Value = 0

$ws.Range("C26").Value = 0
$ws.Range("E26").Value = 0

$ws.Range("C27").Value = 0
$ws.Range("E27").Value = 0

$ws.Range("C28").Value = 0
$ws.Range("E28").Value = 0
